$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# This string appears in the Overview sheet (E2, F2) and in the per-locale
# "Status" column (C2) of the zh-cn and de-de sheets.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- Narrow the Status-related columns to match the new, shorter text ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

# de-de sheet: column C (Status)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
